# Apply the "double file error" fix for Lesson 3 Preparation Answer Key:
#  1. Clarify the bold note applies only to Questions 1-5 (and append a
#     trailing space run, matching the source diff).
#  2. Remove the now-unused "Mode" row from the Problem-1 definitions
#     (Mean / Median / Mode) table block.
#  3. Drop the erroneous " Mode -NN.NNNNNN" figure that had been appended
#     to each Company's Mean/Median solution text.

$d = $word.ActiveDocument

# --- 1. Bold note paragraph -------------------------------------------------
$d.Content.Find.Execute(
    "Please note that the steps show rounded numbers, but that the final answers to the problems are calculated without rounding.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Please note that the steps show rounded numbers, but that the final answers to the problems are calculated without rounding. (Questions 1-5)",
    2)

$notePara = $d.Paragraphs.Item(4)
$noteRange = $notePara.Range
$noteEnd = $d.Range($noteRange.End - 1, $noteRange.End - 1)
$noteEnd.InsertAfter(" ")

# --- 2. Remove the "Mode" definition row (Problem 1 / Part "Mode") ---------
$table = $d.Tables.Item(1)
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    if ($row.Cells.Item(2).Range.Text.TrimEnd([char]7, [char]13) -eq "Mode") {
        $row.Delete()
        break
    }
}

# --- 3. Strip the trailing " Mode -NN.NNNNNN" text from each company row ---
$d.Content.Find.Execute("Mean: 21.276 Median: 13.433 Mode -48.837209", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: 21.276 Median: 13.433", 2)
$d.Content.Find.Execute("Mean: 33.482 Median: 20.838 Mode -62.837689", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: 33.482 Median: 20.838", 2)
$d.Content.Find.Execute("Mean: 41.122 Median: 25.558 Mode -71.050584", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: 41.122 Median: 25.558", 2)
$d.Content.Find.Execute("Mean: 0.706 Median: 1.892 Mode -44.416873", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: 0.706 Median: 1.892", 2)
$d.Content.Find.Execute("Mean: -1.084 Median: -3.796 Mode -39.686099", $true, $false, $false, $false, $false, $true, 1, $false, "Mean: -1.084 Median: -3.796", 2)

Write-Output "done"
